$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Week 2")
$ws3 = $wb.Worksheets.Item("Week 3")

# --- Row 9 (entry 3): fill in date, start, stop, interruption time ---
$ws3.Range("B9").Value = 43511
$ws3.Range("C9").Value = 0.53819444444444442
$ws3.Range("D9").Value = 0.69097222222222221
$ws3.Range("E9").Value = 80

# --- Row 10 (entry 4): fill in date and start time only ---
$ws3.Range("B10").Value = 43512
$ws3.Range("C10").Value = 0.70833333333333337

# Row 8 keeps its wrapped 28.8pt height, now flagged as an explicit
# (author-set) custom height rather than the sheet default
$ws3.Range("A8").RowHeight = 28.8

# --- Copy the Activity/Comments formatting from the matching merged
#     block on "Week 2" so the new merge picks up the same look
#     (top/middle/bottom border treatment around the merged area) ---
$ws2.Range("G7").Copy()
$ws3.Range("G8").PasteSpecial(-4122)

$ws2.Range("H7").Copy()
$ws3.Range("H8").PasteSpecial(-4122)

$ws2.Range("G8").Copy()
$ws3.Range("G9").PasteSpecial(-4122)
$ws3.Range("G10").PasteSpecial(-4122)

$ws2.Range("H8").Copy()
$ws3.Range("H9").PasteSpecial(-4122)
$ws3.Range("H10").PasteSpecial(-4122)

$ws2.Range("G12").Copy()
$ws3.Range("G11").PasteSpecial(-4122)

$ws2.Range("H12").Copy()
$ws3.Range("H11").PasteSpecial(-4122)

# Merge the Activity (G) and Comments (H) columns across rows 8-11,
# matching the shared "Prep." / "Watching JavaScript course..." entry
$ws3.Range("G8:G11").Merge()
$ws3.Range("H8:H11").Merge()

# Leave the new merged block selected, as last touched by the author
$ws3.Range("H8:H11").Select()
